$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellTextPreserveStyle($rng, [string]$val) {
    # Range.Value assignment resets the cell's style (xf) in this
    # engine, so stash the existing format on a scratch cell first,
    # write the new text, then paste the format back.
    $holder = $ws.Range("ZZ1")
    $rng.Copy()
    $holder.PasteSpecial(-4122)
    $rng.Value = $val
    $holder.Copy()
    $rng.PasteSpecial(-4122)
    $holder.Clear()
}

$ws.Range("A2").Value = '01_threequeens_edgecase_flush_str'
$ws.Range("B2").Value = '[(''Q'', ''Spades''), (''Q'', ''Diamonds''), (''2'', ''Diamonds''), (''4'', ''Diamonds''), (''6'', ''Diamonds''), (''8'', ''Diamonds''), (''Q'', ''Hearts''), (''K'', ''Spades''), (''A'', ''Clubs''), (''3'', ''Clubs''), (''6'', ''Hearts''), (''9'', ''Spades''), (''J'', ''Hearts'')]'
$ws.Range("C2").Value = 'Contains 3 Queens where the the Queen of Diamonds is part of a 5 card flush. No fantasy'

$ws.Range("A3").Value = '02_threequeens_edgecase_straight_str'
$ws.Range("B3").Value = '[(''8'', ''Spades''), (''9'', ''Diamonds''), (''10'', ''Clubs''), (''J'', ''Hearts''), (''Q'', ''Hearts''), (''Q'', ''Clubs''), (''Q'', ''Diamonds''), (''10'', ''Hearts''), (''A'', ''Spades''), (''2'', ''Clubs''), (''2'', ''Hearts''), (''3'', ''Diamonds''), (''5'', ''Clubs'')]'
$ws.Range("C3").Value = 'Contains 3 Queens where the Queen rank is part of a 8,9,10,J,Q (Q-high) straight.No Fantasy'

$ws.Range("A4").Value = '03_dragon_hardcode_edgecase'
$ws.Range("B4").Value = '[(''A'', ''Spades''), (''2'', ''Hearts''), (''3'', ''Spades''), (''4'', ''Diamonds''), (''5'', ''Clubs''), (''6'', ''Diamonds''), (''7'', ''Hearts''), (''8'', ''Spades''), (''9'', ''Clubs''), (''10'', ''Clubs''), (''Q'', ''Hearts''), (''K'', ''Spades''), (''J'', ''Hearts'')]'
$ws.Range("C4").Value = 'One card from each rank'

$ws.Range("A5").Value = '04_tripswith2straights_edgecase'
$ws.Range("B5").Value = '[(''10'', ''Diamonds''), (''2'', ''Clubs''),(''3'', ''Hearts''),(''4'', ''Diamonds''),(''5'', ''Hearts''),(''6'', ''Hearts''),(''8'', ''Spades''),(''9'', ''Diamonds''),(''A'', ''Clubs''),(''A'', ''Hearts''),(''J'', ''Spades''),(''A'', ''Diamonds''),(''Q'', ''Diamonds'')]'
$ws.Range("C5").Value = 'A pair and 2,3,4,5,6 and 8,9,10,J,Q straights.Fantasy'

$ws.Range("A6").Value = '05_twoQuads_str'
$ws.Range("B6").Value = '[(''K'', ''Hearts''),(''K'', ''Clubs''),(''9'', ''Hearts''),(''J'', ''Hearts''),(''8'', ''Diamonds''),(''7'', ''Spades''),(''K'', ''Diamonds''),(''7'', ''Hearts''),(''10'', ''Clubs''),(''K'', ''Spades''),(''Q'', ''Hearts''),(''7'', ''Diamonds''),(''7'', ''Clubs'')]'

$ws.Range("A7").Value = '06_2quads and a straight flush'
Set-CellTextPreserveStyle $ws.Range("B7") '[(''10'', ''Hearts''), (''10'', ''Clubs''), (''10'', ''Spades''), (''10'', ''Diamonds''), (''A'', ''Spades''), (''2'', ''Spades''), (''3'', ''Spades''), (''4'', ''Spades''), (''5'', ''Spades''), (''K'', ''Hearts''), (''K'', ''Spades''), (''K'', ''Clubs''), (''K'', ''Diamonds'')]'
$ws.Range("C7").Value = 'triggers fantast for double kings, and quad 10, perhaps develop for all 13 synergy?'

$ws.Range("A8").Value = '07_interaction_flush_got_2_straights_edgecase_str'
Set-CellTextPreserveStyle $ws.Range("B8") '[(''4'', ''Diamonds''), (''5'', ''Hearts''), (''5'', ''Diamonds''), (''6'', ''Spades''), (''7'', ''Spades''), (''7'', ''Hearts''), (''7'', ''Clubs''), (''7'', ''Diamonds''), (''8'', ''Spades''), (''8'', ''Diamonds''), (''9'', ''Diamonds''), (''10'', ''Diamonds''), (''J'', ''Diamonds''),]
'
$ws.Range("C8").Value = 'Contains a 2 straight, and a flush. Classic Case_1 Both Flush and One Long straight from 4 to J inclusive bounds.'

$ws.Range("A9").Value = '08_interaction_flush_got_2_straight_1_good_straight'
Set-CellTextPreserveStyle $ws.Range("B9") '[(''4'', ''Diamonds''), (''5'', ''Hearts''), (''5'', ''Diamonds''), (''6'', ''Spades''), (''7'', ''Spades''), (''7'', ''Hearts''), (''7'', ''Clubs''), (''7'', ''Diamonds''), (''8'', ''Spades''), (''8'', ''Diamonds''), (''9'', ''Hearts''), (''10'', ''Diamonds''), (''J'', ''Diamonds'')]
'
$ws.Range("C9").Value = '(''9'', ''Diamonds'') is subed for (''9'', ''Hearts''). This time there are no loaned_flushed_cards'

$ws.Range("A10").Value = '09_interaction_flush_no2straights_edgecase_str'
Set-CellTextPreserveStyle $ws.Range("B10") '[(''4'', ''Diamonds''), (''5'', ''Hearts''), (''5'', ''Diamonds''), (''6'', ''Spades''), (''6'', ''Diamonds''), (''7'', ''Spades''), (''7'', ''Hearts''), (''7'', ''Clubs''), (''7'', ''Diamonds''), (''8'', ''Spades''), (''9'', ''Hearts''), (''10'', ''Diamonds''), (''J'', ''Diamonds''),]
'
$ws.Range("C10").Value = 'Continueing from prev TC, (''8'', ''Diamonds'') is subed for (''6'', ''Diamonds''). 2 straights no longer found because central cards 7 and 8 are not duplicate. This actually doesn''t factor because it is case_1, expected result bad TC.'

$ws.Range("A11").Value = '10_2_long_straight_involving_ACES'
$ws.Range("B11").Value = '[(''2'', ''Spades'', ), (''3'', ''Diamonds'', ), (''4'', ''Hearts'', ), (''5'', ''Hearts'', ), (''5'', ''Clubs'', ), (''6'', ''Spades'', ), (''9'', ''Clubs'', ), (''10'', ''Hearts'', ), (''J'', ''Diamonds'', ), (''J'', ''Clubs'', ), (''Q'', ''Clubs'', ), (''K'', ''Spades'', ), (''A'', ''Spades'', )]'
$ws.Range("C11").Value = 'A couple of 5''s, Ace-Low, and Ace-High straights. No Fantasy.'

$ws.Range("A12").Value = '11_3_monsters'
$ws.Range("B12").Value = '[(''2'', ''Diamonds''), (''3'', ''Diamonds''), (''6'', ''Diamonds''), (''9'', ''Diamonds''), (''10'', ''Spades''), (''10'', ''Hearts''), (''10'', ''Diamonds''), (''10'', ''Clubs''), (''J'', ''Diamonds''), (''Q'', ''Diamonds''), (''K'', ''Clubs''), (''A'', ''Diamonds''), (''A'', ''Clubs'')]'
$ws.Range("C12").Value = 'Hardcoded 3 monsters, namely quad 10s, diamond flush and 9-A straight.'

$ws.Range("A13").Value = '12_2_flushes_and_2_long_straights'
$ws.Range("B13").Value = '[(''2'', ''Hearts'', ''ID_12''), (''3'', ''Spades'', ''ID_8''), (''3'', ''Diamonds'', ''ID_4''), (''4'', ''Spades'', ''ID_7''), (''5'', ''Diamonds'', ''ID_11''), (''7'', ''Spades'', ''ID_1''), (''8'', ''Diamonds'', ''ID_6''), (''9'', ''Spades'', ''ID_10''), (''10'', ''Diamonds'', ''ID_3''), (''J'', ''Diamonds'', ''ID_9''), (''Q'', ''Spades'', ''ID_2''), (''K'', ''Diamonds'', ''ID_13''), (''A'', ''Diamonds'', ''ID_5'')]'

$ws.Range("B12").Select()

